$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.997.47"
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.305.68"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.55"
$ws.Range("E5").Value = "  +2.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.15"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.48"
$ws.Range("E10").Value = "  +5.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0912"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.979"
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.41"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.657.12"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.303.85"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.899.72"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  +33.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.90"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.61"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "274.24"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.76"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.18"
$ws.Range("E30").Value = "  +11.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.83"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.18"
$ws.Range("E32").Value = "  +7.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0892"
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("E36").Value = "  -11.74%  "
$ws.Range("E37").Value = "  +3.66%  "
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.75"
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("E40").Value = "  -4.50%  "
$ws.Range("E41").Value = "  +8.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.31"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.83"
$ws.Range("E43").Value = "  +3.49%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.226"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.55"
$ws.Range("E46").Value = "  +7.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "83.75"
$ws.Range("E47").Value = "  +12.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "114.56"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.34"
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.92"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.591.71"
$ws.Range("E51").Value = "  +5.24%  "
